$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) — first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 680
$ws1.Range("F3").Value = 520
$ws1.Range("F5").Value = 21
$ws1.Range("F8").Value = 3284
$ws1.Range("F9").Value = 4245
$ws1.Range("F10").Value = 113

# Sheet "全部类型" (All types) — fourth sheet, mirrors the same rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 680
$ws4.Range("F3").Value = 520
$ws4.Range("F5").Value = 21
$ws4.Range("F8").Value = 3284
$ws4.Range("F9").Value = 4245
$ws4.Range("F10").Value = 113
